# Update the single data row (row 2) with the latest report figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE changed from "001" to "002" (REPORT_TYPE_CODE in K2 stays "001").
# Leading apostrophe keeps it a plain text value (no leading-zero loss); reset
# the style afterwards so no stray "quote prefix" number-format is left behind.
$ws.Range("J2").Value = "'002"
$ws.Range("J2").Style = "Normal"

# Report date moved from the 2019 annual report to the 2020 semi-annual report.
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Refreshed balance-sheet figures for the new reporting period.
$ws.Range("O2").Value = 2005451279.16
$ws.Range("P2").Value = 426600116.74
$ws.Range("Q2").Value = 408279917.97
$ws.Range("R2").Value = 97.2631968796
$ws.Range("S2").Value = 341162817.91
$ws.Range("T2").Value = -4.4571479973
$ws.Range("U2").Value = 573173229.11
$ws.Range("V2").Value = -22.4837965661
$ws.Range("W2").Value = 900775596.63
$ws.Range("X2").Value = 177609612.03
$ws.Range("Y2").Value = 22.0595941503

# ADVANCE_RECEIVABLES / ADVANCE_RECEIVABLES_RATIO are not reported this period.
$ws.Range("Z2").Value = ""
$ws.Range("AA2").Value = ""

$ws.Range("AB2").Value = 1104675682.53
$ws.Range("AC2").Value = 10.4600413854
$ws.Range("AD2").Value = -1.5197729873
$ws.Range("AE2").Value = -13.0803858695
$ws.Range("AF2").Value = 159.5333136755
$ws.Range("AG2").Value = 44.9163540391
